$d = $word.ActiveDocument

$old1 = "org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:163)"
$new1 = "org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:162)"
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
Write-Output "Replacement 1 found: $found1"

$old2 = "org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:136)"
$new2 = "org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:135)"
$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
Write-Output "Replacement 2 found: $found2"

$old3 = "org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:168)"
$new3 = "org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:172)"
$found3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
Write-Output "Replacement 3 found: $found3"

$old4 = "org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:118)"
$new4 = "org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:119)"
$found4 = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
Write-Output "Replacement 4 found: $found4"

$old5 = "org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:480)"
$new5 = "org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:462)"
$found5 = $d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)
Write-Output "Replacement 5 found: $found5"

$old6 = "org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:389)"
$new6 = "org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:372)"
$found6 = $d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2)
Write-Output "Replacement 6 found: $found6"

$old7 = "sun.reflect.GeneratedMethodAccessor74.invoke(Unknown Source)"
$new7 = "sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)"
$found7 = $d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1, $false, $new7, 2)
Write-Output "Replacement 7 found: $found7"

$old8 = "`tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:264)`n`tat org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:153)`n`tat org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:124)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:208)`n`tat org.apache.maven.surefire.booter.ProviderFactory`$ProviderProxy.invoke(ProviderFactory.java:156)`n`tat org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:82)`n`tat org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:91)`n`tat org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:587)`n`tat org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:198)`n`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:134)`n`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:104)`n`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:388)`n`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:243)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:656)`n`tat org.eclipse.equinox.launcher.Main.basicRun(Main.java:592)`n`tat org.eclipse.equinox.launcher.Main.run(Main.java:1498)`n`tat org.eclipse.equinox.launcher.Main.main(Main.java:1471)"
$new8 = "`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)`n`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)"
$found8 = $d.Content.Find.Execute($old8, $true, $false, $false, $false, $false, $true, 1, $false, $new8, 2)
Write-Output "Replacement 8 found: $found8"

if (-not ($found1 -and $found2 -and $found3 -and $found4 -and $found5 -and $found6 -and $found7 -and $found8)) {
    throw "One or more expected stack-trace replacements were not found in the document."
}
